$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, pushing the existing row 15 (and below) down to row 16
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with the new record's data
$ws.Cells.Item(15, 1).Value = 4
$ws.Cells.Item(15, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(15, 3).Value = "Los Lagos"
$ws.Cells.Item(15, 4).Value = 44663
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Value = 100112012
$ws.Cells.Item(15, 7).Value = "Espinaca"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 30
$ws.Cells.Item(15, 11).Value = 12000
$ws.Cells.Item(15, 12).Value = 12000
$ws.Cells.Item(15, 13).Value = 12000
$ws.Cells.Item(15, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 1200
$ws.Cells.Item(15, 17).Value = 10
$ws.Cells.Item(15, 18).Value = "Hortaliza"
